# Commit message: "removed ER tags from non-ER templates and non-ER tags"
#
# This RNA-Seq assay template is not an ER (endpoint repository) template,
# so the ER / ER Term Accession Number / ER Term Source REF values on the
# isa_template (formerly "SwateTemplateMetadata") sheet are cleared out,
# leaving the labelled rows empty. The sheet is also renamed to match the
# repo's current naming convention, and becomes the active/selected tab.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("SwateTemplateMetadata")

# Rename the metadata sheet to its new name.
$ws2.Name = "isa_template"

# Clear the ER / ER Term Accession Number / ER Term Source REF values -
# this template doesn't target an endpoint repository.
$ws2.Range("B8").ClearContents()
$ws2.Range("B9").ClearContents()
$ws2.Range("B10").ClearContents()

# The metadata sheet becomes the active tab/selection.
$ws2.Activate()
$ws2.Range("G12").Select()
